$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd email address (shared string + hyperlink display text)
$ws.Range("C8").Value = "aashishss.sachdeva@sitpune.edu.in"
$ws.Hyperlinks.Item(1).TextToDisplay = "aashishss.sachdeva@sitpune.edu.in"

# Add missing role_id value for row 8
$ws.Range("D8").Value = 3

# Make row 7 a custom (taller) height
$ws.Rows.Item(7).RowHeight = 16.4

# Update the selection to row 8 (whole row), active cell A8
$ws.Rows.Item(8).Select()
